$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: updated Fecha (date) values and several price / origin /
# quality fields for Cilantro @ Vega Monumental Concepcion, plus one new
# data row (175) appended at the end of the table.

$ws.Range("D2").Value = 44194
$ws.Range("D3").Value = 44194
$ws.Range("D4").Value = 44336
$ws.Range("D5").Value = 44336
$ws.Range("D6").Value = 44453
$ws.Range("D7").Value = 44453
$ws.Range("D8").Value = 44567
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 700
$ws.Range("M8").Value = 650
$ws.Range("N8").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O8").Value = 'Región de Ñuble'
$ws.Range("P8").Value = 650
$ws.Range("Q8").Value = 1
$ws.Range("D9").Value = 44567
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = 500
$ws.Range("P9").Value = 500
$ws.Range("D10").Value = 44308
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 700
$ws.Range("M10").Value = 650
$ws.Range("P10").Value = 650
$ws.Range("D11").Value = 44308
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("M11").Value = 500
$ws.Range("P11").Value = 500
$ws.Range("D12").Value = 44398
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = 650
$ws.Range("P12").Value = 650
$ws.Range("D13").Value = 44398
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = 500
$ws.Range("P13").Value = 500
$ws.Range("D14").Value = 44278
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 600
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 650
$ws.Range("P14").Value = 650
$ws.Range("D15").Value = 44278
$ws.Range("I15").Value = 'Segunda'
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 500
$ws.Range("P15").Value = 500
$ws.Range("D16").Value = 44376
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = 650
$ws.Range("P16").Value = 650
$ws.Range("D17").Value = 44376
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = 500
$ws.Range("P17").Value = 500
$ws.Range("D18").Value = 44574
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 600
$ws.Range("L18").Value = 700
$ws.Range("M18").Value = 650
$ws.Range("P18").Value = 650
$ws.Range("D19").Value = 44574
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = 500
$ws.Range("P19").Value = 500
$ws.Range("D20").Value = 44204
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 700
$ws.Range("M20").Value = 650
$ws.Range("P20").Value = 650
$ws.Range("D21").Value = 44204
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 500
$ws.Range("M21").Value = 500
$ws.Range("P21").Value = 500
$ws.Range("D22").Value = 44320
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = 650
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 650
$ws.Range("D23").Value = 44320
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 500
$ws.Range("L23").Value = 500
$ws.Range("M23").Value = 500
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 500
$ws.Range("D24").Value = 44657
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 180
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 6500
$ws.Range("M24").Value = 6222
$ws.Range("N24").Value = '$/caja 36 atados'
$ws.Range("O24").Value = 'Región Metropolitana'
$ws.Range("P24").Value = 173
$ws.Range("Q24").Value = 36
$ws.Range("D25").Value = 44264
$ws.Range("D26").Value = 44264
$ws.Range("D27").Value = 44168
$ws.Range("O27").Value = 'Región de Ñuble'
$ws.Range("D28").Value = 44168
$ws.Range("O28").Value = 'Región de Ñuble'
$ws.Range("D29").Value = 44391
$ws.Range("O29").Value = 'Región de Ñuble'
$ws.Range("D30").Value = 44391
$ws.Range("O30").Value = 'Región de Ñuble'
$ws.Range("D31").Value = 44509
$ws.Range("D32").Value = 44509
$ws.Range("D33").Value = 44460
$ws.Range("D34").Value = 44460
$ws.Range("D35").Value = 44344
$ws.Range("D36").Value = 44344
$ws.Range("D37").Value = 44614
$ws.Range("D38").Value = 44614
$ws.Range("D39").Value = 44306
$ws.Range("D40").Value = 44306
$ws.Range("D41").Value = 44313
$ws.Range("D42").Value = 44313
$ws.Range("D43").Value = 44327
$ws.Range("D44").Value = 44327
$ws.Range("D45").Value = 44525
$ws.Range("D46").Value = 44525
$ws.Range("D47").Value = 44490
$ws.Range("D48").Value = 44490
$ws.Range("D49").Value = 44636
$ws.Range("J49").Value = 80
$ws.Range("K49").Value = 8000
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = 8375
$ws.Range("N49").Value = '$/caja 36 atados'
$ws.Range("O49").Value = 'Región Metropolitana'
$ws.Range("P49").Value = 233
$ws.Range("Q49").Value = 36
$ws.Range("D50").Value = 44322
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 600
$ws.Range("L50").Value = 700
$ws.Range("M50").Value = 650
$ws.Range("P50").Value = 650
$ws.Range("D51").Value = 44322
$ws.Range("I51").Value = 'Segunda'
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 500
$ws.Range("M51").Value = 500
$ws.Range("P51").Value = 500
$ws.Range("D52").Value = 44425
$ws.Range("I52").Value = 'Primera'
$ws.Range("J52").Value = 200
$ws.Range("K52").Value = 600
$ws.Range("L52").Value = 700
$ws.Range("M52").Value = 650
$ws.Range("P52").Value = 650
$ws.Range("D53").Value = 44425
$ws.Range("I53").Value = 'Segunda'
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 500
$ws.Range("L53").Value = 500
$ws.Range("M53").Value = 500
$ws.Range("P53").Value = 500
$ws.Range("D54").Value = 44579
$ws.Range("I54").Value = 'Primera'
$ws.Range("J54").Value = 200
$ws.Range("K54").Value = 600
$ws.Range("L54").Value = 700
$ws.Range("M54").Value = 650
$ws.Range("O54").Value = 'Región Metropolitana'
$ws.Range("P54").Value = 650
$ws.Range("D55").Value = 44579
$ws.Range("I55").Value = 'Segunda'
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = 500
$ws.Range("O55").Value = 'Región Metropolitana'
$ws.Range("P55").Value = 500
$ws.Range("D56").Value = 44166
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 200
$ws.Range("K56").Value = 600
$ws.Range("L56").Value = 700
$ws.Range("M56").Value = 650
$ws.Range("P56").Value = 650
$ws.Range("D57").Value = 44166
$ws.Range("I57").Value = 'Segunda'
$ws.Range("J57").Value = 100
$ws.Range("K57").Value = 500
$ws.Range("L57").Value = 500
$ws.Range("M57").Value = 500
$ws.Range("P57").Value = 500
$ws.Range("D58").Value = 44433
$ws.Range("I58").Value = 'Primera'
$ws.Range("J58").Value = 200
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 700
$ws.Range("M58").Value = 650
$ws.Range("P58").Value = 650
$ws.Range("D59").Value = 44433
$ws.Range("I59").Value = 'Segunda'
$ws.Range("J59").Value = 100
$ws.Range("K59").Value = 500
$ws.Range("L59").Value = 500
$ws.Range("M59").Value = 500
$ws.Range("P59").Value = 500
$ws.Range("D60").Value = 44209
$ws.Range("I60").Value = 'Primera'
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 700
$ws.Range("M60").Value = 650
$ws.Range("P60").Value = 650
$ws.Range("D61").Value = 44209
$ws.Range("I61").Value = 'Segunda'
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 500
$ws.Range("L61").Value = 500
$ws.Range("M61").Value = 500
$ws.Range("P61").Value = 500
$ws.Range("D62").Value = 44350
$ws.Range("I62").Value = 'Primera'
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 700
$ws.Range("M62").Value = 650
$ws.Range("P62").Value = 650
$ws.Range("D63").Value = 44350
$ws.Range("I63").Value = 'Segunda'
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 500
$ws.Range("L63").Value = 500
$ws.Range("M63").Value = 500
$ws.Range("P63").Value = 500
$ws.Range("D64").Value = 44476
$ws.Range("I64").Value = 'Primera'
$ws.Range("J64").Value = 200
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 700
$ws.Range("M64").Value = 650
$ws.Range("P64").Value = 650
$ws.Range("D65").Value = 44476
$ws.Range("I65").Value = 'Segunda'
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 500
$ws.Range("L65").Value = 500
$ws.Range("M65").Value = 500
$ws.Range("P65").Value = 500
$ws.Range("D66").Value = 44565
$ws.Range("I66").Value = 'Primera'
$ws.Range("J66").Value = 200
$ws.Range("K66").Value = 600
$ws.Range("L66").Value = 700
$ws.Range("M66").Value = 650
$ws.Range("P66").Value = 650
$ws.Range("D67").Value = 44565
$ws.Range("I67").Value = 'Segunda'
$ws.Range("J67").Value = 100
$ws.Range("K67").Value = 500
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 500
$ws.Range("P67").Value = 500
$ws.Range("D68").Value = 44161
$ws.Range("I68").Value = 'Primera'
$ws.Range("J68").Value = 200
$ws.Range("K68").Value = 600
$ws.Range("L68").Value = 700
$ws.Range("M68").Value = 650
$ws.Range("P68").Value = 650
$ws.Range("D69").Value = 44161
$ws.Range("I69").Value = 'Segunda'
$ws.Range("J69").Value = 100
$ws.Range("K69").Value = 500
$ws.Range("L69").Value = 500
$ws.Range("M69").Value = 500
$ws.Range("P69").Value = 500
$ws.Range("D70").Value = 44316
$ws.Range("I70").Value = 'Primera'
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 600
$ws.Range("L70").Value = 700
$ws.Range("M70").Value = 650
$ws.Range("O70").Value = 'Región Metropolitana'
$ws.Range("P70").Value = 650
$ws.Range("D71").Value = 44316
$ws.Range("I71").Value = 'Segunda'
$ws.Range("J71").Value = 100
$ws.Range("K71").Value = 500
$ws.Range("L71").Value = 500
$ws.Range("M71").Value = 500
$ws.Range("O71").Value = 'Región Metropolitana'
$ws.Range("P71").Value = 500
$ws.Range("D72").Value = 44420
$ws.Range("I72").Value = 'Primera'
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 600
$ws.Range("L72").Value = 700
$ws.Range("M72").Value = 650
$ws.Range("P72").Value = 650
$ws.Range("D73").Value = 44420
$ws.Range("I73").Value = 'Segunda'
$ws.Range("J73").Value = 100
$ws.Range("K73").Value = 500
$ws.Range("L73").Value = 500
$ws.Range("M73").Value = 500
$ws.Range("P73").Value = 500
$ws.Range("D74").Value = 44169
$ws.Range("I74").Value = 'Primera'
$ws.Range("J74").Value = 200
$ws.Range("K74").Value = 600
$ws.Range("L74").Value = 700
$ws.Range("M74").Value = 650
$ws.Range("P74").Value = 650
$ws.Range("D75").Value = 44169
$ws.Range("I75").Value = 'Segunda'
$ws.Range("J75").Value = 100
$ws.Range("K75").Value = 500
$ws.Range("L75").Value = 500
$ws.Range("M75").Value = 500
$ws.Range("P75").Value = 500
$ws.Range("D76").Value = 44245
$ws.Range("I76").Value = 'Primera'
$ws.Range("J76").Value = 200
$ws.Range("K76").Value = 600
$ws.Range("L76").Value = 700
$ws.Range("M76").Value = 650
$ws.Range("P76").Value = 650
$ws.Range("D77").Value = 44245
$ws.Range("I77").Value = 'Segunda'
$ws.Range("J77").Value = 100
$ws.Range("K77").Value = 500
$ws.Range("L77").Value = 500
$ws.Range("M77").Value = 500
$ws.Range("P77").Value = 500
$ws.Range("D78").Value = 44656
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 130
$ws.Range("K78").Value = 650
$ws.Range("L78").Value = 6000
$ws.Range("M78").Value = 3942
$ws.Range("N78").Value = '$/caja 36 atados'
$ws.Range("O78").Value = 'Región Metropolitana'
$ws.Range("P78").Value = 110
$ws.Range("Q78").Value = 36
$ws.Range("D79").Value = 44488
$ws.Range("J79").Value = 200
$ws.Range("K79").Value = 600
$ws.Range("L79").Value = 700
$ws.Range("M79").Value = 650
$ws.Range("N79").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O79").Value = 'Región de Ñuble'
$ws.Range("P79").Value = 650
$ws.Range("Q79").Value = 1
$ws.Range("D80").Value = 44488
$ws.Range("I80").Value = 'Segunda'
$ws.Range("J80").Value = 100
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 500
$ws.Range("M80").Value = 500
$ws.Range("P80").Value = 500
$ws.Range("D81").Value = 44250
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 200
$ws.Range("K81").Value = 600
$ws.Range("L81").Value = 700
$ws.Range("M81").Value = 650
$ws.Range("P81").Value = 650
$ws.Range("D82").Value = 44250
$ws.Range("I82").Value = 'Segunda'
$ws.Range("J82").Value = 100
$ws.Range("K82").Value = 500
$ws.Range("L82").Value = 500
$ws.Range("M82").Value = 500
$ws.Range("P82").Value = 500
$ws.Range("D83").Value = 44222
$ws.Range("I83").Value = 'Primera'
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 700
$ws.Range("M83").Value = 650
$ws.Range("P83").Value = 650
$ws.Range("D84").Value = 44222
$ws.Range("I84").Value = 'Segunda'
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 500
$ws.Range("L84").Value = 500
$ws.Range("M84").Value = 500
$ws.Range("P84").Value = 500
$ws.Range("D85").Value = 44239
$ws.Range("I85").Value = 'Primera'
$ws.Range("J85").Value = 200
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 700
$ws.Range("M85").Value = 650
$ws.Range("P85").Value = 650
$ws.Range("D86").Value = 44239
$ws.Range("I86").Value = 'Segunda'
$ws.Range("J86").Value = 100
$ws.Range("K86").Value = 500
$ws.Range("L86").Value = 500
$ws.Range("M86").Value = 500
$ws.Range("P86").Value = 500
$ws.Range("D87").Value = 44217
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 200
$ws.Range("K87").Value = 600
$ws.Range("L87").Value = 700
$ws.Range("M87").Value = 650
$ws.Range("P87").Value = 650
$ws.Range("D88").Value = 44217
$ws.Range("I88").Value = 'Segunda'
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 500
$ws.Range("L88").Value = 500
$ws.Range("M88").Value = 500
$ws.Range("P88").Value = 500
$ws.Range("D89").Value = 44469
$ws.Range("I89").Value = 'Primera'
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 600
$ws.Range("L89").Value = 700
$ws.Range("M89").Value = 650
$ws.Range("P89").Value = 650
$ws.Range("D90").Value = 44469
$ws.Range("I90").Value = 'Segunda'
$ws.Range("J90").Value = 100
$ws.Range("K90").Value = 500
$ws.Range("L90").Value = 500
$ws.Range("M90").Value = 500
$ws.Range("P90").Value = 500
$ws.Range("D91").Value = 44267
$ws.Range("I91").Value = 'Primera'
$ws.Range("J91").Value = 200
$ws.Range("K91").Value = 600
$ws.Range("L91").Value = 700
$ws.Range("M91").Value = 650
$ws.Range("P91").Value = 650
$ws.Range("D92").Value = 44267
$ws.Range("I92").Value = 'Segunda'
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 500
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 500
$ws.Range("P92").Value = 500
$ws.Range("D93").Value = 44334
$ws.Range("I93").Value = 'Primera'
$ws.Range("J93").Value = 200
$ws.Range("K93").Value = 600
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = 650
$ws.Range("P93").Value = 650
$ws.Range("D94").Value = 44334
$ws.Range("I94").Value = 'Segunda'
$ws.Range("J94").Value = 100
$ws.Range("K94").Value = 500
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = 500
$ws.Range("P94").Value = 500
$ws.Range("D95").Value = 44475
$ws.Range("I95").Value = 'Primera'
$ws.Range("J95").Value = 200
$ws.Range("K95").Value = 600
$ws.Range("L95").Value = 700
$ws.Range("M95").Value = 650
$ws.Range("P95").Value = 650
$ws.Range("D96").Value = 44475
$ws.Range("I96").Value = 'Segunda'
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 500
$ws.Range("L96").Value = 500
$ws.Range("M96").Value = 500
$ws.Range("P96").Value = 500
$ws.Range("D97").Value = 44523
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = 650
$ws.Range("P97").Value = 650
$ws.Range("D98").Value = 44523
$ws.Range("I98").Value = 'Segunda'
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 500
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 500
$ws.Range("O98").Value = 'Región de Ñuble'
$ws.Range("P98").Value = 500
$ws.Range("D99").Value = 44607
$ws.Range("I99").Value = 'Primera'
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 700
$ws.Range("M99").Value = 650
$ws.Range("O99").Value = 'Región de Ñuble'
$ws.Range("P99").Value = 650
$ws.Range("D100").Value = 44607
$ws.Range("I100").Value = 'Segunda'
$ws.Range("J100").Value = 100
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 500
$ws.Range("M100").Value = 500
$ws.Range("P100").Value = 500
$ws.Range("D101").Value = 44595
$ws.Range("I101").Value = 'Primera'
$ws.Range("J101").Value = 200
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 700
$ws.Range("M101").Value = 650
$ws.Range("P101").Value = 650
$ws.Range("D102").Value = 44595
$ws.Range("I102").Value = 'Segunda'
$ws.Range("J102").Value = 100
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 500
$ws.Range("M102").Value = 500
$ws.Range("P102").Value = 500
$ws.Range("D103").Value = 44442
$ws.Range("I103").Value = 'Primera'
$ws.Range("J103").Value = 200
$ws.Range("K103").Value = 600
$ws.Range("L103").Value = 700
$ws.Range("M103").Value = 650
$ws.Range("P103").Value = 650
$ws.Range("D104").Value = 44442
$ws.Range("I104").Value = 'Segunda'
$ws.Range("J104").Value = 100
$ws.Range("K104").Value = 500
$ws.Range("L104").Value = 500
$ws.Range("M104").Value = 500
$ws.Range("P104").Value = 500
$ws.Range("D105").Value = 44462
$ws.Range("I105").Value = 'Primera'
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 600
$ws.Range("L105").Value = 700
$ws.Range("M105").Value = 650
$ws.Range("P105").Value = 650
$ws.Range("D106").Value = 44462
$ws.Range("I106").Value = 'Segunda'
$ws.Range("J106").Value = 100
$ws.Range("K106").Value = 500
$ws.Range("L106").Value = 500
$ws.Range("M106").Value = 500
$ws.Range("P106").Value = 500
$ws.Range("D107").Value = 44292
$ws.Range("I107").Value = 'Primera'
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 650
$ws.Range("P107").Value = 650
$ws.Range("D108").Value = 44292
$ws.Range("I108").Value = 'Segunda'
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 500
$ws.Range("L108").Value = 500
$ws.Range("M108").Value = 500
$ws.Range("P108").Value = 500
$ws.Range("D109").Value = 44159
$ws.Range("I109").Value = 'Primera'
$ws.Range("J109").Value = 200
$ws.Range("K109").Value = 600
$ws.Range("L109").Value = 700
$ws.Range("M109").Value = 650
$ws.Range("P109").Value = 650
$ws.Range("D110").Value = 44159
$ws.Range("I110").Value = 'Segunda'
$ws.Range("J110").Value = 100
$ws.Range("K110").Value = 500
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 500
$ws.Range("P110").Value = 500
$ws.Range("D111").Value = 44274
$ws.Range("I111").Value = 'Primera'
$ws.Range("J111").Value = 200
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 700
$ws.Range("M111").Value = 650
$ws.Range("P111").Value = 650
$ws.Range("D112").Value = 44274
$ws.Range("I112").Value = 'Segunda'
$ws.Range("J112").Value = 100
$ws.Range("K112").Value = 500
$ws.Range("L112").Value = 500
$ws.Range("M112").Value = 500
$ws.Range("P112").Value = 500
$ws.Range("D113").Value = 44187
$ws.Range("I113").Value = 'Primera'
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 650
$ws.Range("P113").Value = 650
$ws.Range("D114").Value = 44187
$ws.Range("I114").Value = 'Segunda'
$ws.Range("J114").Value = 100
$ws.Range("K114").Value = 500
$ws.Range("L114").Value = 500
$ws.Range("M114").Value = 500
$ws.Range("P114").Value = 500
$ws.Range("D115").Value = 44645
$ws.Range("I115").Value = 'Primera'
$ws.Range("J115").Value = 140
$ws.Range("K115").Value = 6500
$ws.Range("L115").Value = 7000
$ws.Range("M115").Value = 6714
$ws.Range("N115").Value = '$/caja 36 atados'
$ws.Range("O115").Value = 'Región Metropolitana'
$ws.Range("P115").Value = 186
$ws.Range("Q115").Value = 36
$ws.Range("D116").Value = 44405
$ws.Range("D117").Value = 44405
$ws.Range("D118").Value = 44365
$ws.Range("D119").Value = 44365
$ws.Range("D120").Value = 44341
$ws.Range("D121").Value = 44341
$ws.Range("D122").Value = 44280
$ws.Range("D123").Value = 44280
$ws.Range("D124").Value = 44362
$ws.Range("D125").Value = 44362
$ws.Range("D126").Value = 44285
$ws.Range("D127").Value = 44285
$ws.Range("D128").Value = 44271
$ws.Range("D129").Value = 44271
$ws.Range("D130").Value = 44512
$ws.Range("D131").Value = 44512
$ws.Range("D132").Value = 44582
$ws.Range("D133").Value = 44582
$ws.Range("D134").Value = 44447
$ws.Range("D135").Value = 44447
$ws.Range("D136").Value = 44299
$ws.Range("D137").Value = 44299
$ws.Range("D138").Value = 44434
$ws.Range("D139").Value = 44434
$ws.Range("D140").Value = 44371
$ws.Range("D141").Value = 44371
$ws.Range("D142").Value = 44266
$ws.Range("D143").Value = 44266
$ws.Range("D144").Value = 44231
$ws.Range("D145").Value = 44231
$ws.Range("D146").Value = 44237
$ws.Range("D147").Value = 44237
$ws.Range("D148").Value = 44224
$ws.Range("D149").Value = 44224
$ws.Range("D150").Value = 44553
$ws.Range("D151").Value = 44553
$ws.Range("D152").Value = 44616
$ws.Range("D153").Value = 44616
$ws.Range("D154").Value = 44386
$ws.Range("D155").Value = 44386
$ws.Range("D156").Value = 44427
$ws.Range("D157").Value = 44427
$ws.Range("D158").Value = 44252
$ws.Range("D159").Value = 44252
$ws.Range("D160").Value = 44330
$ws.Range("D161").Value = 44330
$ws.Range("D162").Value = 44609
$ws.Range("D163").Value = 44609
$ws.Range("D164").Value = 44203
$ws.Range("D165").Value = 44203
$ws.Range("D166").Value = 44435
$ws.Range("J166").Value = 400
$ws.Range("D167").Value = 44435
$ws.Range("J167").Value = 200
$ws.Range("D168").Value = 44358
$ws.Range("D169").Value = 44358
$ws.Range("D170").Value = 44467
$ws.Range("D171").Value = 44467
$ws.Range("D172").Value = 44257
$ws.Range("D173").Value = 44257
$ws.Range("D174").Value = 44383
$ws.Range("J174").Value = 200
$ws.Range("K174").Value = 600
$ws.Range("L174").Value = 700
$ws.Range("M174").Value = 650
$ws.Range("N174").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O174").Value = 'Región de Ñuble'
$ws.Range("P174").Value = 650
$ws.Range("Q174").Value = 1
$ws.Range("A175").Value = 11
$ws.Range("B175").Value = 'Vega Monumental Concepción'
$ws.Range("C175").Value = 'Bíobío'
$ws.Range("D175").Value = 44383
$ws.Range("E175").Value = 8
$ws.Range("F175").Value = 100112040
$ws.Range("G175").Value = 'Cilantro'
$ws.Range("H175").Value = 'Sin especificar'
$ws.Range("I175").Value = 'Segunda'
$ws.Range("J175").Value = 100
$ws.Range("K175").Value = 500
$ws.Range("L175").Value = 500
$ws.Range("M175").Value = 500
$ws.Range("N175").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O175").Value = 'Región de Ñuble'
$ws.Range("P175").Value = 500
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = 'Hortaliza'
$ws.Range("D175").NumberFormat = "YYYY-MM-DD HH:MM:SS"
